$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column keeps its existing text formatting so values
# like "0.6930" or "1.000" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.093.56"
$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("D3").Value = "1.848.92"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "0.6930"
$ws.Range("E5").Value = "  -4.99%  "
$ws.Range("D6").Value = "237.90"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.07722"
$ws.Range("E8").Value = "  +8.89%  "
$ws.Range("D9").Value = "0.3031"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "23.22"
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("D11").Value = "0.08107"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "0.7244"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").Value = "1.839.58"
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("D14").Value = "5.205"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "88.94"
$ws.Range("E15").Value = "  -3.31%  "
$ws.Range("D16").Value = "29.086.81"
$ws.Range("E16").Value = "  -2.30%  "
$ws.Range("D17").Value = "5.742"
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").Value = "0.000007764"
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "13.18"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "235.79"
$ws.Range("E20").Value = "  -4.77%  "
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "2.093.67"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D24").Value = "7.580"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "8.969"
$ws.Range("E25").Value = "  -2.18%  "
$ws.Range("D26").Value = "160.87"
$ws.Range("E26").Value = "  -1.28%  "
$ws.Range("D27").Value = "0.1430"
$ws.Range("E27").Value = "  -7.08%  "
$ws.Range("D28").Value = "18.02"
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("D29").Value = "1.980"
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "1.398"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("D31").Value = "4.488"
$ws.Range("E31").Value = "  -0.70%  "
$ws.Range("D32").Value = "1.484"
$ws.Range("E32").Value = "  -2.33%  "
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("D34").Value = "0.05220"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("D35").Value = "1.183"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").Value = "0.6992"
$ws.Range("E37").Value = "  -5.92%  "
$ws.Range("D38").Value = "2.656"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").Value = "0.01851"
$ws.Range("E39").Value = "  -3.94%  "
$ws.Range("D40").Value = "2.679"
$ws.Range("E40").Value = "  -1.80%  "
$ws.Range("D41").Value = "0.9182"
$ws.Range("E41").Value = "  +6.23%  "
$ws.Range("D42").Value = "6.008"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").Value = "1.079.81"
$ws.Range("E43").Value = "  +3.60%  "
$ws.Range("D44").Value = "0.4249"
$ws.Range("E44").Value = "  -4.44%  "
$ws.Range("D45").Value = "70.26"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "103.11"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "1.768"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").Value = "1.990.87"
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("D50").Value = "9.121"
$ws.Range("E50").Value = "  -3.94%  "
$ws.Range("D51").Value = "6.991"
$ws.Range("E51").Value = "  -6.03%  "
